# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.970.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.414.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "639.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.400"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.964"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.414.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("E12").Value = "  -4.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.973.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.055.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000249"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -10.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.421.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "498.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.468"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.33%  "

$ws.Range("E24").Value = "  -5.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000191"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "90.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.600.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.87%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.79%  "

$ws.Range("E33").Value = "  -6.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.176"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.89%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.65%  "

$ws.Range("E37").Value = "  -3.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "539.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.62%  "

$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  -1.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.903"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.55%  "

$ws.Range("E45").Value = "  +0.91%  "

$ws.Range("E46").Value = "  -1.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0403"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "
